$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

$ws.Range("D2").Value = "67.364.05"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.624.64"
$ws.Range("E3").Value = "  -1.49%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.01%  "
Set-TextValue $ws.Range("D5") "593.59"
$ws.Range("E5").Value = "  -0.63%  "
Set-TextValue $ws.Range("D6") "168.01"
$ws.Range("E6").Value = "  +1.60%  "
$ws.Range("E7").Value = "  +0.03%  "
Set-TextValue $ws.Range("D8") "0.533"
$ws.Range("E8").Value = "  -2.42%  "
$ws.Range("D9").Value = "2.623.63"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("E11").Value = "  +1.12%  "
Set-TextValue $ws.Range("D12") "0.364"
$ws.Range("E12").Value = "  +2.05%  "
Set-TextValue $ws.Range("D13") "5.23"
$ws.Range("E13").Value = "  +0.27%  "
Set-TextValue $ws.Range("D14") "27.64"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "3.106.50"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "67.162.54"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "2.629.48"
$ws.Range("E18").Value = "  -0.49%  "
Set-TextValue $ws.Range("D19") "12.05"
$ws.Range("E19").Value = "  +3.10%  "
Set-TextValue $ws.Range("D20") "7.97"
$ws.Range("E20").Value = "  +4.59%  "
Set-TextValue $ws.Range("D21") "355.90"
$ws.Range("E21").Value = "  -1.89%  "
Set-TextValue $ws.Range("D22") "4.31"
$ws.Range("E22").Value = "  -1.24%  "
Set-TextValue $ws.Range("D23") "4.66"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("E24").Value = "  -0.04%  "
Set-TextValue $ws.Range("D25") "1.93"
$ws.Range("E25").Value = "  -4.87%  "
Set-TextValue $ws.Range("D26") "10.24"
$ws.Range("E26").Value = "  +2.12%  "
Set-TextValue $ws.Range("D27") "69.61"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").Value = "2.761.27"
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("E30").Value = "  -1.48%  "
Set-TextValue $ws.Range("D31") "544.16"
$ws.Range("E31").Value = "  -1.93%  "
Set-TextValue $ws.Range("D32") "7.91"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("E33").Value = "  -2.52%  "
Set-TextValue $ws.Range("D34") "1.90"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("E35").Value = "  +4.97%  "
$ws.Range("E36").Value = "  +0.03%  "
Set-TextValue $ws.Range("D37") "1.50"
$ws.Range("E37").Value = "  -3.11%  "
Set-TextValue $ws.Range("D38") "156.24"
$ws.Range("E38").Value = "  +0.10%  "
Set-TextValue $ws.Range("D39") "19.01"
$ws.Range("E39").Value = "  -2.40%  "
Set-TextValue $ws.Range("D40") "0.366"
$ws.Range("E40").Value = "  -1.74%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D41") "1.81"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D42") "18.20"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue $ws.Range("D43") "5.21"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("E44").Value = "  +0.08%  "
Set-TextValue $ws.Range("D45") "2.42"
$ws.Range("E45").Value = "  -3.51%  "
$ws.Range("D46").Value = "0.0₆0298"
$ws.Range("E46").Value = "  +0.36%  "
Set-TextValue $ws.Range("D47") "152.69"
$ws.Range("E47").Value = "  +0.05%  "
Set-TextValue $ws.Range("D48") "0.580"
$ws.Range("E48").Value = "  -1.32%  "
Set-TextValue $ws.Range("D49") "3.78"
$ws.Range("E49").Value = "  -1.06%  "
Set-TextValue $ws.Range("D50") "1.70"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("E51").Value = "  -1.27%  "
